$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1800
$ws.Range("I100").Value = 1787.5
$ws.Range("J100").Value = 1850
$ws.Range("K100").Value = 1787.5
$ws.Range("L100").Value = 1850
$ws.Range("M100").Value = -1246.5
$ws.Range("N100").Value = -2932

$ws.Range("H107").Value = 513.9167
$ws.Range("I107").Value = 406.1
$ws.Range("J107").Value = 1053
$ws.Range("K107").Value = 406.1
$ws.Range("L107").Value = 1053
$ws.Range("M107").Value = 1513.9
$ws.Range("N107").Value = -4893

$ws.Range("H112").Value = 37038850
$ws.Range("J112").Value = 2009.9565
$ws.Range("L112").Value = 6029.8695
$ws.Range("N112").Value = -8245.8695

$ws.Range("H129").Value = 997
$ws.Range("I129").Value = 662.8571
$ws.Range("J129").Value = 1113.95
$ws.Range("K129").Value = 1988.5713
$ws.Range("L129").Value = 3341.85
$ws.Range("M129").Value = 3011.4287
$ws.Range("N129").Value = -13341.85

$ws.Range("H132").Value = 767200.56
$ws.Range("I132").Value = 1359.0358
$ws.Range("K132").Value = 4077.1074
$ws.Range("M132").Value = -1547.1074

$ws.Range("H137").Value = 2501625.2
$ws.Range("I137").Value = 4546803
$ws.Range("J137").Value = 1963.6666
$ws.Range("K137").Value = 13640409
$ws.Range("L137").Value = 5890.9998
$ws.Range("M137").Value = -13637859
$ws.Range("N137").Value = -10990.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2721.4
$ws.Range("J2").Value = 2599.8572
$ws.Range("L2").Value = 2599.8572
$ws.Range("N2").Value = -2825.8572

$ws.Range("H97").Value = 2404678
$ws.Range("I97").Value = 3677289.8
$ws.Range("J97").Value = 855.55554
$ws.Range("K97").Value = 3677289.8
$ws.Range("L97").Value = 855.55554
$ws.Range("M97").Value = -3676793.8
$ws.Range("N97").Value = -1847.55554

$ws.Range("H102").Value = 6499753.5
$ws.Range("I102").Value = 7525556.5
$ws.Range("K102").Value = 7525556.5
$ws.Range("M102").Value = -7523934.5

$ws.Range("H110").Value = 910570.56
$ws.Range("I110").Value = 2001096
$ws.Range("K110").Value = 2001096
$ws.Range("M110").Value = -1999051

$ws.Range("H116").Value = 2721.4
$ws.Range("J116").Value = 2599.8572
$ws.Range("L116").Value = 2599.8572
$ws.Range("N116").Value = -7187.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2721.4
$ws.Range("J3").Value = 2599.8572
$ws.Range("L3").Value = 2599.8572
$ws.Range("N3").Value = -2827.8572

$ws.Range("H40").Value = 24500
$ws.Range("J40").Value = 24500
$ws.Range("L40").Value = 24500
$ws.Range("N40").Value = -25030

$ws.Range("H107").Value = 1784.4474
$ws.Range("I107").Value = 1755.15
$ws.Range("J107").Value = 1817
$ws.Range("K107").Value = 1755.15
$ws.Range("L107").Value = 1817
$ws.Range("M107").Value = 164.8499999999999
$ws.Range("N107").Value = -5657

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3201.2046
$ws.Range("I31").Value = 1420.3914
$ws.Range("J31").Value = 5151.619
$ws.Range("K31").Value = 1420.3914
$ws.Range("L31").Value = 5151.619
$ws.Range("M31").Value = -1125.3914
$ws.Range("N31").Value = -5741.619

$ws.Range("H34").Value = 3201.2046
$ws.Range("I34").Value = 1420.3914
$ws.Range("J34").Value = 5151.619
$ws.Range("K34").Value = 1420.3914
$ws.Range("L34").Value = 5151.619
$ws.Range("M34").Value = -1218.3914
$ws.Range("N34").Value = -5555.619

$ws.Range("H107").Value = 603.6667
$ws.Range("I107").Value = 603.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 603.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1316.3333
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 18413.809
$ws.Range("I134").Value = 1292.6888
$ws.Range("J134").Value = 61216.61
$ws.Range("K134").Value = 3878.0664
$ws.Range("L134").Value = 183649.83
$ws.Range("M134").Value = -1343.0664
$ws.Range("N134").Value = -188719.83

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2597.75
$ws.Range("I76").Value = 91
$ws.Range("J76").Value = 3433.3333
$ws.Range("K76").Value = 273
$ws.Range("L76").Value = 10299.9999
$ws.Range("M76").Value = 110
$ws.Range("N76").Value = -11065.9999

$ws.Range("H79").Value = 2597.75
$ws.Range("I79").Value = 91
$ws.Range("J79").Value = 3433.3333
$ws.Range("K79").Value = 273
$ws.Range("L79").Value = 10299.9999
$ws.Range("M79").Value = 1053
$ws.Range("N79").Value = -12951.9999

$ws.Range("H80").Value = 1930.3914
$ws.Range("I80").Value = 1199.8572
$ws.Range("K80").Value = 3599.5716
$ws.Range("M80").Value = -2663.5716

$ws.Range("H82").Value = 3600
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3600
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 10800
$ws.Range("N82").Value = -11612
$ws.Range("M82").ClearContents()

$ws.Range("H83").Value = 1930.3914
$ws.Range("I83").Value = 1199.8572
$ws.Range("K83").Value = 10798.7148
$ws.Range("M83").Value = -6118.7148

$ws.Range("H85").Value = 3600
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3600
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 10800
$ws.Range("N85").Value = -13608
$ws.Range("M85").ClearContents()

$ws.Range("H86").Value = 1580
$ws.Range("J86").Value = 2300
$ws.Range("L86").Value = 6900
$ws.Range("N86").Value = -9272

$ws.Range("H87").Value = 27373.943
$ws.Range("I87").Value = 17669.166
$ws.Range("J87").Value = 29381.828
$ws.Range("K87").Value = 53007.49800000001
$ws.Range("L87").Value = 88145.484
$ws.Range("M87").Value = -51759.49800000001
$ws.Range("N87").Value = -90641.484

$ws.Range("H88").Value = 2607.6924
$ws.Range("I88").Value = 500
$ws.Range("J88").Value = 2990.9092
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 8972.7276
$ws.Range("M88").Value = -1072
$ws.Range("N88").Value = -9828.7276

$ws.Range("H89").Value = 1580
$ws.Range("J89").Value = 2300
$ws.Range("L89").Value = 20700
$ws.Range("N89").Value = -32556

$ws.Range("H90").Value = 27373.943
$ws.Range("I90").Value = 17669.166
$ws.Range("J90").Value = 29381.828
$ws.Range("K90").Value = 159022.494
$ws.Range("L90").Value = 264436.452
$ws.Range("M90").Value = -152782.494
$ws.Range("N90").Value = -276916.452

$ws.Range("H91").Value = 2607.6924
$ws.Range("I91").Value = 500
$ws.Range("J91").Value = 2990.9092
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 8972.7276
$ws.Range("M91").Value = -18
$ws.Range("N91").Value = -11936.7276

$ws.Range("H131").Value = 894
$ws.Range("J131").Value = 1072.3055
$ws.Range("L131").Value = 3216.9165
$ws.Range("N131").Value = -13296.9165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3464.875
$ws.Range("I80").Value = 2783.75
$ws.Range("K80").Value = 2783.75
$ws.Range("M80").Value = -1785.75

$ws.Range("H83").Value = 3464.875
$ws.Range("I83").Value = 2783.75
$ws.Range("K83").Value = 13918.75
$ws.Range("M83").Value = -8926.75

$ws.Range("H97").Value = 2363.75
$ws.Range("I97").Value = 2363.75
$ws.Range("K97").Value = 2363.75
$ws.Range("M97").Value = -1867.75

$ws.Range("H107").Value = 428.45456
$ws.Range("I107").Value = 210
$ws.Range("J107").Value = 510.375
$ws.Range("K107").Value = 210
$ws.Range("L107").Value = 510.375
$ws.Range("M107").Value = 1710
$ws.Range("N107").Value = -4350.375

$ws.Range("H113").Value = 1933.3334
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 2080
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 2080
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -6420

$ws.Range("H123").Value = 24483.25
$ws.Range("J123").Value = 24483.25
$ws.Range("L123").Value = 24483.25
$ws.Range("N123").Value = -29383.25

$ws.Range("H132").Value = 51428.977
$ws.Range("I132").Value = 33702.934
$ws.Range("J132").Value = 112485.336
$ws.Range("K132").Value = 101108.802
$ws.Range("L132").Value = 337456.008
$ws.Range("M132").Value = -98578.802
$ws.Range("N132").Value = -342516.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3751.2083
$ws.Range("I16").Value = 1135.2354
$ws.Range("K16").Value = 1135.2354
$ws.Range("M16").Value = -965.2354

$ws.Range("H22").Value = 964.4167
$ws.Range("I22").Value = 467.27274
$ws.Range("J22").Value = 1385.0769
$ws.Range("K22").Value = 467.27274
$ws.Range("L22").Value = 1385.0769
$ws.Range("M22").Value = -172.27274
$ws.Range("N22").Value = -1975.0769

$ws.Range("H27").Value = 964.4167
$ws.Range("I27").Value = 467.27274
$ws.Range("J27").Value = 1385.0769
$ws.Range("K27").Value = 467.27274
$ws.Range("L27").Value = 1385.0769
$ws.Range("M27").Value = -360.27274
$ws.Range("N27").Value = -1599.0769

$ws.Range("H61").Value = 2668.7646
$ws.Range("I61").Value = 2520.6924
$ws.Range("J61").Value = 3150
$ws.Range("K61").Value = 2520.6924
$ws.Range("L61").Value = 3150
$ws.Range("M61").Value = -2318.6924
$ws.Range("N61").Value = -3554

$ws.Range("H113").Value = 2668.7646
$ws.Range("I113").Value = 2520.6924
$ws.Range("J113").Value = 3150
$ws.Range("K113").Value = 2520.6924
$ws.Range("L113").Value = 3150
$ws.Range("M113").Value = -350.6923999999999
$ws.Range("N113").Value = -7490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 24715.4
$ws.Range("J119").Value = 24715.4
$ws.Range("L119").Value = 24715.4
$ws.Range("N119").Value = -34391.4

$ws.Range("H132").Value = 52416.15
$ws.Range("I132").Value = 39217.367
$ws.Range("J132").Value = 138208.25
$ws.Range("K132").Value = 117652.101
$ws.Range("L132").Value = 414624.75
$ws.Range("M132").Value = -115122.101
$ws.Range("N132").Value = -419684.75
